$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1): new columns L1:X1 need the border/bold style ---
# clone style from an existing styled header cell (B1) into the new header cells
foreach ($col in @("L","M","N","O","P","Q","R","S","T","U","V","W","X")) {
    $ws.Range("B1").Copy($ws.Range($col + "1"))
}

# --- Set header text for every header cell B1:X1 ---
$headers = [ordered]@{
    "B1" = '(''N'', ''H'', 1)'
    "C1" = '(''N'', ''C'', 1)'
    "D1" = '(''C'', ''N'', 3)'
    "E1" = '(''C'', ''S'', 1)'
    "F1" = '(''O'', ''C'', 1.5)'
    "G1" = '(''F'', ''C'', 1)'
    "H1" = '(''C'', ''O'', 2)'
    "I1" = '(''C'', ''C'', 1)'
    "J1" = '(''C'', ''C'', 1.5)'
    "K1" = '(''C'', ''H'', 1)'
    "L1" = 'C'
    "M1" = 'F'
    "N1" = 'H'
    "O1" = 'N'
    "P1" = 'O'
    "Q1" = 'amines'
    "R1" = 'aldehyde'
    "S1" = 'esters'
    "T1" = 'cyanide'
    "U1" = 'aromatic'
    "V1" = 'weight'
    "W1" = 'logP'
    "X1" = 'SAS'
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- Set data rows 2-5 across columns A-X ---
$rows = [ordered]@{
    2 = @(0, 1.222, 5.632, 0.068, 0.065, 0.159, 0.344, 1.821, 12.618, 5.806, 25.543, 20.337, 0.344, 26.93, 2.352, 3.049, 1.958, 1.821, 0.113, 0.068, 0.968, 360.344, 2.636, 3.322)
    3 = @(1, 1.23, 0.36, 0.153, 0.681, 0.307, 1.03, 0.879, 1.867, 6.473, 4.456, 9.695, 1.03, 6.261, 2.021, 2.301, 0.346, 0.879, 0.193, 0.153, 1.736, 270.905, 2.364, 2.416)
    4 = @(2, 2.424, 8.275, 0.087, 0.191, 0.141, 0.698, 1.514, 6.602, 1.537, 14.53, 12.215, 0.698, 17.15, 4.421, 1.384, 2.75, 1.514, 0.052, 0.087, 0.482, 291, 1.154, 3.082)
    5 = @(3, 0.151, 1.555, 0.098, 0.382, 0.193, 0.717, 1.872, 4.231, 14.767, 15.633, 19.922, 0.717, 16.026, 0.794, 4.755, 0.846, 1.872, 0.212, 0.098, 2.552, 381.678, 3.584, 1.919)
}
foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $colLetter = [string]([char](65 + $i))
        $ws.Range($colLetter + $r).Value = $vals[$i]
    }
}

# --- Remove old row 6 (table now only has rows 1-5) ---
$ws.Rows.Item(6).Delete()
